$d = $word.ActiveDocument

# --- Need section: merge '.  It ' + 'could ' runs, drop gramStart proofErr ---
$d.Content.Find.Execute(".  It could ", $true, $false, $false, $false, $false, $true, 1, $false, ".  It could ", 1) | Out-Null

# --- Need section: merge 'be used' + ' ' runs, drop gramEnd proofErr ---
$d.Content.Find.Execute("be used ", $true, $false, $false, $false, $false, $true, 1, $false, "be used ", 1) | Out-Null

# --- Approach section: merge 'is based' sentence, drop gramStart/gramEnd proofErr ---
$d.Content.Find.Execute("In order to meet this need we plan to create an anonymous communication network that is based upon the popular IRC protocol.", $true, $false, $false, $false, $false, $true, 1, $false, "In order to meet this need we plan to create an anonymous communication network that is based upon the popular IRC protocol.", 1) | Out-Null

# --- Approach section: merge 'can be identified' sentence, drop gramStart/gramEnd proofErr ---
$d.Content.Find.Execute("  Since people can be identified in this method, we will need to research ways to circumvent this.", $true, $false, $false, $false, $false, $true, 1, $false, "  Since people can be identified in this method, we will need to research ways to circumvent this.", 1) | Out-Null

# --- Competition section: rewrite paragraph with updated content about Anonychat ---
$old5 = "The idea of an anonymous communication system is not anything new, and there have been numerous attempts to implement one. The Invisible Internet Project (I2P) contains an anonymous IRC protocol; however, I2P focuses more on an overall anonymous communication systems, whereas we are just focusing on a standalone communication program. Freenode is another anonymous IRC system, however Freenode still makes use of a centralized server, and uses SSL encryption rather than distribution for anonymity. A few others exists as well, such as Quassel and Rust, however these projects are no longer under active development and have taken different approaches."
$new5 = "The idea of an anonymous communication network is not new, however our Anonychat will still bring together a combination of things that does not currently exist on other products. Some anonymous communications include the Invisible Internet Project (I2P), an ongoing effort to build a free, open source, and anonymous internet. I2P includes a system to allow anonymous IRC communication, by simply allowing standard IRC protocol over the I2P network. Since I2P is designed at the network layer, it does restrict compatibility to only those on the I2P network to maintain anonymity, whereas Anonychat’s restrictions are to the application itself. Users will be able connect across any existing communication network. Freenode is an example of IRC using the standard protocol with SSL encryption to ensure anonymity. Freenode, being an IRC protocol, still will have direct connections that can indicate relations between users, while our Anonychat will implement a method to obfuscate intended targets of messages. In addition, Anonychat aims to be more directly peer to peer, requiring a central server only for initial connections to the Anonychat network. Competiton also extends to peer to peer style communication network Skype, which uses a similar connection system we intend to implement (a central server to start, then p2p communication afterwards). A key difference between Skype and Anonychat is that Skype’s main focus is not anonymity, and uses P2P connections in a more direct method. Some projects in the IRC field, such as Quassel and Rust, are more direct competition in anonymity, but do not use the P2P connectivity we will attempt. Overall, our metric for success will be if Anonychat is able to reliably send messages through the distributed system with minimal chance for the messages to directly connect users."
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

Write-Output "done"
